$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Busqueda" in E1, reusing the same header formatting (style)
# applied to the other header cells (A1:D1).
$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "Busqueda"

# Update the date/time in column A for rows 2-11, and add the "percy" search
# term used to find these products in the new column E.
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = "2025-05-22 17:56"
    $ws.Cells.Item($r, 5).Value = "percy"
}
